$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "/resolve"
$ws.Range("A9").Value = "/duplicate-from-path"
$ws.Range("B10").Value = "/resolve"
$ws.Range("A10").Value = "/duplicate-from-path"

$ws.Range("B11").Select()

$ws.Columns.Item(1).AutoFit() | Out-Null

$wb.Save()
